$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.240143656730652
$ws.Range("B1").Value = 2.357472658157349
$ws.Range("C1").Value = 3.201707601547241
$ws.Range("D1").Value = 3.544786930084229
$ws.Range("E1").Value = 1.107582211494446
